# Separate environment URL data out of credData.xlsx.
# The five "*URL" columns (BOURL, ofbizURL, CCURL, CWURL, RSURL) are being
# moved out to a separate environment-data file, so here we simply delete
# those columns from Sheet1, leaving only the user/pass pairs behind.
# Deleting from right to left keeps the remaining column letters stable
# while each delete is performed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("M").Delete()   # RSURL
$ws.Columns("J").Delete()   # CWURL
$ws.Columns("G").Delete()   # CCURL
$ws.Columns("D").Delete()   # ofbizURL
$ws.Columns("A").Delete()   # BOURL

# The worksheet's hyperlinks collection still references the old (now
# stale) cell locations -- rebuild it from scratch so only the two
# surviving hyperlinked cells (the C@bi$ush5 mailto + the michigan@na.com
# mailto) keep their links, now anchored at their new positions.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:C@bi`$ush5", "", "", "C@bi`$ush5") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "mailto:michigan@na.com") | Out-Null

# Adding a hyperlink with a display string overwrites the cell text with
# that display string, and re-applies the "Hyperlink" cell style as a
# brand-new style record -- restore the original cell text ("cabiautomation")
# for D2 (the hyperlink keeps its own separate "display" attribute) and put
# back the plain "Hyperlink" style on both cells so they keep using the same
# cell style as before (same font, no border) instead of a stray duplicate.
$ws.Range("D2").Value2 = "cabiautomation"
$ws.Range("D2").Style = "Hyperlink"
$ws.Range("I2").Style = "Hyperlink"

# Match the post-edit selection left behind in the saved file (column I
# fully selected, no frozen/scrolled top-left cell).
$ws.Columns("I").Select()
